# Rename the inline picture shapes' docPr "name" attributes, matching the
# source-controlled OOXML rename:
#   footer (both the "odd/primary" and "even" footer) Pearson logo:
#       image1.png -> image2.png
#   header (the header that actually carries the BTec logo) :
#       image2.jpg -> image1.jpg
#
# Pictures live as inline drawings inside the section's headers/footers, so
# we reach them via Section.Headers / Section.Footers -> Range.InlineShapes.

$d = $word.ActiveDocument
$sec = $d.Sections.Item(1)

# --- Footers: both footers hold the Pearson logo (image1.png -> image2.png) ---
for ($f = 1; $f -le $sec.Footers.Count; $f++) {
    $ftr = $sec.Footers.Item($f)
    if ($ftr.Exists) {
        $shapes = $ftr.Range.InlineShapes
        for ($i = 1; $i -le $shapes.Count; $i++) {
            $shp = $shapes.Item($i)
            if ($shp.AlternativeText -eq "Y:\Together Design\Pearson Edexcel PowerPoint amends\Assets\PearsonLogo.png") {
                $shp.Name = "image2.png"
            }
        }
    }
}

# --- Headers: the header with the BTec logo (image2.jpg -> image1.jpg) ---
for ($h = 1; $h -le $sec.Headers.Count; $h++) {
    $hdr = $sec.Headers.Item($h)
    if ($hdr.Exists) {
        $shapes = $hdr.Range.InlineShapes
        for ($i = 1; $i -le $shapes.Count; $i++) {
            $shp = $shapes.Item($i)
            if ($shp.AlternativeText -eq "BTec_Logo-Orange") {
                $shp.Name = "image1.jpg"
            }
        }
    }
}
